# Apply targeted data corrections to the "Admin Guide Export" sheet.
# Rows refer directly to worksheet rows (row 1 = header).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin Guide Export")

# Row 7 - "Analyze App Metadata Analyzer": clear the redundant Tooling value.
$ws.Range("H7").Value = $null

# Row 17 - "Audit User Access": simplify Tooling Options to just "Native".
$ws.Range("G17").Value = "Native"

# Row 34 - "Review Pinning/Load Balancing": simplify Tooling Options and
# clear the redundant Tooling value.
$ws.Range("G34").Value = "Native"
$ws.Range("H34").Value = $null

# Row 36 - "Review/Update Capacity Plan": simplify Tooling Options and
# clear the redundant Tooling value.
$ws.Range("G36").Value = "Native"
$ws.Range("H36").Value = $null

# Update the active selection to match the saved view state.
$ws.Range("G32").Select()
